$d = $word.ActiveDocument

$d.Content.Find.Execute(
    ">>>  your stuff after this line >>>",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "This is testing for version control. I will do the pull request now",
    2
)
